$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: set a full row of H..N values (column 8..14) on a worksheet.
# Any value passed as $null means "leave the cell empty / clear it".
# NOTE: uses positional parameters only (named parameter binding for
# user-defined functions is not reliable in this runtime).
# ---------------------------------------------------------------------------
function Set-RowHN {
    param($ws, [int]$Row, $H, $I, $J, $K, $L, $M, $N)

    $values = @($H, $I, $J, $K, $L, $M, $N)
    for ($i = 0; $i -lt 7; $i++) {
        $col = 8 + $i            # H=8 .. N=14
        $val = $values[$i]
        $cell = $ws.Cells.Item($Row, $col)
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# ===========================================================================
# Sheet "ALC" (sheet1.xml)
# ===========================================================================
$wsALC = $wb.Worksheets.Item("ALC")

# Row 116: update H/I/K, drop M
Set-RowHN $wsALC 116 2355.0908 0 2355.0908 0 2355.0908 $null -9239.0908

# Row 137: update H..N
Set-RowHN $wsALC 137 839.5833 702.6316 1360 2107.8948 4080 442.1052 -9180

# ===========================================================================
# Sheet "BSM" (sheet3.xml)
# ===========================================================================
$wsBSM = $wb.Worksheets.Item("BSM")

Set-RowHN $wsBSM 117 30000     0 30000     0 30000     $null -39178
Set-RowHN $wsBSM 118 39000     0 39000     0 39000     $null -42314
Set-RowHN $wsBSM 119 29000     0 29000     0 29000     $null -38676
Set-RowHN $wsBSM 120 27904     0 27904     0 27904     $null -37580
Set-RowHN $wsBSM 122 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 123 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 124 48000     0 48000     0 48000     $null -57820
Set-RowHN $wsBSM 125 65000     0 65000     0 65000     $null -74840
Set-RowHN $wsBSM 126 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 127 58000     0 58000     0 58000     $null -67920
Set-RowHN $wsBSM 128 1574.75   1574.75 0   4724.25 0   -2234.25 $null
Set-RowHN $wsBSM 129 49949.25  0 49949.25  0 49949.25  $null -59949.25
Set-RowHN $wsBSM 130 33714.145 0 33714.145 0 33714.145 $null -43754.145
Set-RowHN $wsBSM 131 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 132 62096.668 0 62096.668 0 62096.668 $null -72216.66800000001
Set-RowHN $wsBSM 133 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 134 2680.0715 2825.0908 2148.3333 8475.2724 6444.999899999999 -5940.2724 -11514.9999
Set-RowHN $wsBSM 135 49725     0 49725     0 49725     $null -59865
Set-RowHN $wsBSM 137 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 138 78564.28999999999 0 78564.28999999999 0 78564.28999999999 $null -88844.28999999999
Set-RowHN $wsBSM 139 0         0 0         0 0         $null $null
Set-RowHN $wsBSM 140 60000     0 60000     0 60000     $null -70360
Set-RowHN $wsBSM 141 42445     0 42445     0 42445     $null -52805

# ===========================================================================
# Sheet "CRP" (sheet4.xml)
# ===========================================================================
$wsCRP = $wb.Worksheets.Item("CRP")

# Rows 31 and 34: update H..N
Set-RowHN $wsCRP 31 1051.53 1301.2142 870.7241 1301.2142 870.7241 -1006.2142 -1460.7241
Set-RowHN $wsCRP 34 1051.53 1301.2142 870.7241 1301.2142 870.7241 -1099.2142 -1274.7241

Set-RowHN $wsCRP 129 43185.285 0 43185.285 0 43185.285 $null -53185.285
Set-RowHN $wsCRP 130 74780.71000000001 0 74780.71000000001 0 74780.71000000001 $null -84820.71000000001
Set-RowHN $wsCRP 131 47394.5 0 47394.5 0 47394.5 $null -57474.5
Set-RowHN $wsCRP 132 2099.9355 1488.4615 5279.6 4465.3845 15838.8 -1935.3845 -20898.8
Set-RowHN $wsCRP 133 37000 0 37000 0 37000 $null -42060
Set-RowHN $wsCRP 134 1044.2667 1013.6667 1166.6666 3041.0001 3499.9998 -506.0001000000002 -8569.9998
Set-RowHN $wsCRP 135 0 0 0 0 0 $null $null
Set-RowHN $wsCRP 137 40780 0 40780 0 40780 $null -50980
Set-RowHN $wsCRP 138 0 0 0 0 0 $null $null
Set-RowHN $wsCRP 139 0 0 0 0 0 $null $null
Set-RowHN $wsCRP 140 35000 0 35000 0 35000 $null -45360
Set-RowHN $wsCRP 141 55744.5 0 55744.5 0 55744.5 $null -66104.5

# ===========================================================================
# Sheet "CUL" (sheet5.xml)
# ===========================================================================
$wsCUL = $wb.Worksheets.Item("CUL")

Set-RowHN $wsCUL 131 4313007.5 438.57144 5901848.5 1315.71432 17705545.5 3724.28568 -17715625.5

# ===========================================================================
# Sheet "LTW" (sheet7.xml)
# ===========================================================================
$wsLTW = $wb.Worksheets.Item("LTW")

Set-RowHN $wsLTW 124 63500 0 63500 0 63500 $null -73320
Set-RowHN $wsLTW 125 0 0 0 0 0 $null $null
Set-RowHN $wsLTW 127 0 0 0 0 0 $null $null
Set-RowHN $wsLTW 128 40000 0 40000 0 40000 $null -49960
Set-RowHN $wsLTW 129 61050 0 61050 0 61050 $null -71050
Set-RowHN $wsLTW 130 47500 0 47500 0 47500 $null -57540
Set-RowHN $wsLTW 131 60000 60000 0 60000 0 -54960 $null
Set-RowHN $wsLTW 132 11230.04 15340.2 5064.8 46020.60000000001 15194.4 -43490.60000000001 -20254.4
Set-RowHN $wsLTW 133 0 0 0 0 0 $null $null
Set-RowHN $wsLTW 134 40000 0 40000 0 40000 $null -50140
Set-RowHN $wsLTW 135 50000 0 50000 0 50000 $null -60140
Set-RowHN $wsLTW 136 12399.5 26500 2999.1667 79500 8997.500100000001 -76950 -14097.5001
Set-RowHN $wsLTW 137 47429 0 47429 0 47429 $null -57629
Set-RowHN $wsLTW 138 0 0 0 0 0 $null $null
Set-RowHN $wsLTW 139 45589.375 30000 47816.43 30000 47816.43 -24860 -58096.43
Set-RowHN $wsLTW 140 35000 0 35000 0 35000 $null -45360
Set-RowHN $wsLTW 141 79950 0 79950 0 79950 $null -90310
